{"js": "// Remove the \"Requisitos\" Heading2 paragraph and the following\n// \"LOT2040 -  Engenharia Gen\u00e9tica  (Requisito fraco)\" ListBullet\n// paragraph that follows the Bibliografia section, per the diff:\n// both paragraphs are deleted entirely, leaving the Bibliografia\n// paragraph as the last paragraph in the document body.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text,items/style\");\nawait context.sync();\n\nconst toDelete = [];\nfor (const p of paragraphs.items) {\n  const text = (p.text || \"\").trim();\n  if (\n    text === \"Requisitos\" ||\n    text.indexOf(\"LOT2040 -  Engenharia Gen\u00e9tica  (Requisito fraco)\") === 0\n  ) {\n    toDelete.push(p);\n  }\n}\n\nfor (const p of toDelete) {\n  p.delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the \"Requisitos\" Heading2 paragraph and the following\n# \"LOT2040 -  Engenharia Gen\u00e9tica  (Requisito fraco)\" ListBullet\n# paragraph (the final two paragraphs of the document, right after\n# the Bibliografia section), per the diff.\n\n$d = $word.ActiveDocument\n\nfor ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {\n    $p = $d.Paragraphs.Item($i)\n    $text = $p.Range.Text.Trim()\n    if ($text -eq \"Requisitos\" -or $text.StartsWith(\"LOT2040 -  Engenharia Gen\u00e9tica  (Requisito fraco)\")) {\n        $p.Range.Delete()\n    }\n}\n"}
